$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")
$ws.Rows.Item(56).Resize(4).Delete() | Out-Null
Write-Host ($ws.Range("A56").Value())
